# Update the marksheet's correct/total marks figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row: total points awarded for right answers (B11): 3 -> 5
$ws.Range("B11").Value = 5

# "Total" row: total score (B12): 51 -> 85
$ws.Range("B12").Value = 85

# "Total" row: Correct/Max display text (E12): "49/84" -> "85/140"
$ws.Range("E12").Value = "85/140"
